# Applies the "fpga design calcs.xlsx" edit:
#  - Sheet1: re-enter a few formulas (previously Excel auto-"shared"), add new
#    L79:L82 PRODUCT block, add v/dn & dn/a labeled values (M85/N85, M86/N86),
#    add the "3MHZ 4CH 12BIT..." string in M91 and a hex-dump fill (L92:M141).
#  - Divide: re-enter B29/C29/P29/Q29 (drop their "shared" grouping).
#  - Model Alignment: re-enter C14/C15/C16 (drop their "shared" grouping).
#  - Sheet1 view: scroll & selection moved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# Re-enter these formulas plainly (author re-typed / re-filled them instead of
# letting Excel's drag-fill keep them grouped as "shared" formulas).
$ws1.Range("L32").Formula = "=POWER(2,K32)*0.2005/(205*3000000)*100"

$ws1.Range("L39").Formula = "=POWER(2,K39)*200.5/205/3"
$ws1.Range("K40").Formula = "=K39+1"
$ws1.Range("L40").Formula = "=POWER(2,K40)*200.5/205/3"
$ws1.Range("K41").Formula = "=K40+1"
$ws1.Range("L41").Formula = "=POWER(2,K41)*200.5/205/3"
$ws1.Range("K42").Formula = "=K41+1"
$ws1.Range("L42").Formula = "=POWER(2,K42)*200.5/205/3"

# New small PRODUCT block next to the existing F79:F82 calc.
$ws1.Range("L79").Value = 6
$ws1.Range("L80").Value = 8
$ws1.Range("L81").Value = 256
$ws1.Range("L82").Formula = "=PRODUCT(L79:L81)"

# New annotated copies of F84/F85 over in M/N with labels.
$ws1.Range("M85").Value = 0.2005
$ws1.Range("N85").Value = "v/dn"
$ws1.Range("M86").Value = 205
$ws1.Range("N86").Value = "dn/a"

# New string + hex-dump-by-character fill below the existing log calc.
$ws1.Range("M91").Value = "3MHZ 4CH 12BIT 4MSample TRACE BUFFER, 800x480 XVGA"

$ws1.Range("L92").Value = 1
$ws1.Range("M92").Formula = '=DEC2HEX(CODE(MID($M$91,L92,1)),2)'
for ($row = 93; $row -le 141; $row++) {
    $prevRow = $row - 1
    $ws1.Range("L$row").Formula = "=L$prevRow+1"
    $ws1.Range("M$row").Formula = "=DEC2HEX(CODE(MID(`$M`$91,L$row,1)),2)"
}

# View moved while scrolled down working on the above.
$ws1.Range("I146").Select()

# ---------------------------------------------------------------------------
# Divide
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Divide")
$ws2.Range("B29").Formula = "=B28+1"
$ws2.Range("C29").Formula = "=C28-2"
$ws2.Range("P29").Formula = "=P28-2"
$ws2.Range("Q29").Formula = "=Q28-2"

# ---------------------------------------------------------------------------
# Model Alignment
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Model Alignment")
$ws3.Range("C14").Formula = "=SUM(E14:BF14)"
$ws3.Range("C15").Formula = "=SUM(E15:BF15)"
$ws3.Range("C16").Formula = "=SUM(E16:BF16)"

# Re-select Sheet1 as the active sheet (it was last active before save).
$ws1.Activate()
